# Apply "updated data from quant engine" refresh to the PSU Fund holdings comparison sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: month columns roll forward (Jan_2026 added, Oct_2025 dropped) ---
$ws.Cells.Item(1,4).Value = "Jan_2026"
$ws.Cells.Item(1,5).Value = "Dec_2025"
$ws.Cells.Item(1,6).Value = "Nov_2025"
# Columns G (MoM) and H (QoQ) headers are unchanged.

# --- Data rows: refreshed holdings, 22 rows (was 24) ---
$ws.Cells.Item(2,1).Value = "INE522F01014"
$ws.Cells.Item(2,2).Value = "Coal India Ltd"
$ws.Cells.Item(2,3).Value = "quant PSU Fund"
$ws.Cells.Item(2,4).Value = 10.101646
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 10.101646
$ws.Cells.Item(2,8).Value = 10.101646

$ws.Cells.Item(3,1).Value = "INE115A01026"
$ws.Cells.Item(3,2).Value = "LIC Housing Finance Ltd"
$ws.Cells.Item(3,3).Value = "quant PSU Fund"
$ws.Cells.Item(3,4).Value = 9.69026
$ws.Cells.Item(3,5).Value = 9.588673
$ws.Cells.Item(3,6).Value = 9.196228
$ws.Cells.Item(3,7).Value = 0.1015870000000003
$ws.Cells.Item(3,8).Value = 0.4940320000000007

$ws.Cells.Item(4,1).Value = "INE0J1Y01017"
$ws.Cells.Item(4,2).Value = "Life Insurance Corporation Of India"
$ws.Cells.Item(4,3).Value = "quant PSU Fund"
$ws.Cells.Item(4,4).Value = 9.38495
$ws.Cells.Item(4,5).Value = 8.925459
$ws.Cells.Item(4,6).Value = 8.800407999999999
$ws.Cells.Item(4,7).Value = 0.4594909999999999
$ws.Cells.Item(4,8).Value = 0.5845420000000008

$ws.Cells.Item(5,1).Value = "INE423A01024"
$ws.Cells.Item(5,2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(5,3).Value = "quant PSU Fund"
$ws.Cells.Item(5,4).Value = 8.559979999999999
$ws.Cells.Item(5,5).Value = 8.703604
$ws.Cells.Item(5,6).Value = 6.043951
$ws.Cells.Item(5,7).Value = -0.1436240000000009
$ws.Cells.Item(5,8).Value = 2.516029

$ws.Cells.Item(6,1).Value = "INE364U01010"
$ws.Cells.Item(6,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(6,3).Value = "quant PSU Fund"
$ws.Cells.Item(6,4).Value = 7.921467
$ws.Cells.Item(6,5).Value = 8.653581000000001
$ws.Cells.Item(6,6).Value = 3.631247
$ws.Cells.Item(6,7).Value = -0.732114000000001
$ws.Cells.Item(6,8).Value = 4.29022

$ws.Cells.Item(7,1).Value = "INE131A01031"
$ws.Cells.Item(7,2).Value = "Gujarat Mineral Development Corp Ltd"
$ws.Cells.Item(7,3).Value = "quant PSU Fund"
$ws.Cells.Item(7,4).Value = 7.343367
$ws.Cells.Item(7,5).Value = 6.98031
$ws.Cells.Item(7,6).Value = 5.925037
$ws.Cells.Item(7,7).Value = 0.3630569999999995
$ws.Cells.Item(7,8).Value = 1.41833

$ws.Cells.Item(8,1).Value = "INE584A01023"
$ws.Cells.Item(8,2).Value = "NMDC Ltd"
$ws.Cells.Item(8,3).Value = "quant PSU Fund"
$ws.Cells.Item(8,4).Value = 7.217659
$ws.Cells.Item(8,5).Value = 6.779965
$ws.Cells.Item(8,6).Value = 2.87107
$ws.Cells.Item(8,7).Value = 0.4376940000000005
$ws.Cells.Item(8,8).Value = 4.346589

$ws.Cells.Item(9,1).Value = "INE213A01029"
$ws.Cells.Item(9,2).Value = "Oil and Natural Gas Corporation Ltd."
$ws.Cells.Item(9,3).Value = "quant PSU Fund"
$ws.Cells.Item(9,4).Value = 6.520862
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 6.520862
$ws.Cells.Item(9,8).Value = 6.520862

$ws.Cells.Item(10,1).Value = "INE018E01016"
$ws.Cells.Item(10,2).Value = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(10,3).Value = "quant PSU Fund"
$ws.Cells.Item(10,4).Value = 2.931496
$ws.Cells.Item(10,5).Value = 3.074732
$ws.Cells.Item(10,6).Value = 1.517926
$ws.Cells.Item(10,7).Value = -0.1432359999999999
$ws.Cells.Item(10,8).Value = 1.41357

$ws.Cells.Item(11,1).Value = "INE589A01014"
$ws.Cells.Item(11,2).Value = "NLC India Limited"
$ws.Cells.Item(11,3).Value = "quant PSU Fund"
$ws.Cells.Item(11,4).Value = 0.603333
$ws.Cells.Item(11,5).Value = 0.534424
$ws.Cells.Item(11,6).Value = 0.489507
$ws.Cells.Item(11,7).Value = 0.068909
$ws.Cells.Item(11,8).Value = 0.113826

$ws.Cells.Item(12,1).Value = "INE752E01010"
$ws.Cells.Item(12,2).Value = "Power Grid Corporation of India Limited"
$ws.Cells.Item(12,3).Value = "quant PSU Fund"
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,6).Value = 2.711656
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = -2.711656

$ws.Cells.Item(13,1).Value = "INE510A01028"
$ws.Cells.Item(13,2).Value = "Engineers India Limited"
$ws.Cells.Item(13,3).Value = "quant PSU Fund"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 1.011007
$ws.Cells.Item(13,6).Value = 3.066512
$ws.Cells.Item(13,7).Value = -1.011007
$ws.Cells.Item(13,8).Value = -3.066512

$ws.Cells.Item(14,1).Value = "INE423A20016"
$ws.Cells.Item(14,2).Value = "Adani Enterprises Limited Rights"
$ws.Cells.Item(14,3).Value = "quant PSU Fund"
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 0.14071
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = -0.14071

$ws.Cells.Item(15,1).Value = "INE257A01026"
$ws.Cells.Item(15,2).Value = "Bharat Heavy Electricals Ltd"
$ws.Cells.Item(15,3).Value = "quant PSU Fund"
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 8.221399999999999
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = -8.221399999999999
$ws.Cells.Item(15,8).Value = 0

$ws.Cells.Item(16,1).Value = "INE263A01024"
$ws.Cells.Item(16,2).Value = "Bharat Electronics Ltd"
$ws.Cells.Item(16,3).Value = "quant PSU Fund"
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 9.592911000000001
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = -9.592911000000001

$ws.Cells.Item(17,1).Value = "INE020B01018"
$ws.Cells.Item(17,2).Value = "Rural Electrification Corporation Ltd"
$ws.Cells.Item(17,3).Value = "quant PSU Fund"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 0.75745
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = -0.75745

$ws.Cells.Item(18,1).Value = "INE171Z01026"
$ws.Cells.Item(18,2).Value = "Bharat Dynamics Limited"
$ws.Cells.Item(18,3).Value = "quant PSU Fund"
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 3.631098
$ws.Cells.Item(18,6).Value = 1.952966
$ws.Cells.Item(18,7).Value = -3.631098
$ws.Cells.Item(18,8).Value = -1.952966

$ws.Cells.Item(19,1).Value = "INE134E01011"
$ws.Cells.Item(19,2).Value = "Power Finance Corporation Ltd."
$ws.Cells.Item(19,3).Value = "quant PSU Fund"
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 2.951911
$ws.Cells.Item(19,7).Value = 0
$ws.Cells.Item(19,8).Value = -2.951911

$ws.Cells.Item(20,1).Value = "INE095N01031"
$ws.Cells.Item(20,2).Value = "National Building Construction Corp"
$ws.Cells.Item(20,3).Value = "quant PSU Fund"
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 1.515132
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(20,8).Value = -1.515132

$ws.Cells.Item(21,1).Value = "INE062A01020"
$ws.Cells.Item(21,2).Value = "State Bank of India"
$ws.Cells.Item(21,3).Value = "quant PSU Fund"
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 9.260284
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = -9.260284
$ws.Cells.Item(21,8).Value = 0

$ws.Cells.Item(22,1).Value = "INE031A01017"
$ws.Cells.Item(22,2).Value = "Housing & Urban Devlopment Company Ltd"
$ws.Cells.Item(22,3).Value = "quant PSU Fund"
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 6.050591
$ws.Cells.Item(22,6).Value = 5.980718
$ws.Cells.Item(22,7).Value = -6.050591
$ws.Cells.Item(22,8).Value = -5.980718

$ws.Cells.Item(23,1).Value = "INE931S01010"
$ws.Cells.Item(23,2).Value = "Adani Energy Solutions Limited"
$ws.Cells.Item(23,3).Value = "quant PSU Fund"
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 1.826353
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = -1.826353

# The old sheet had 24 data rows (through row 25); the refreshed data only has 22
# (through row 23), so remove the two now-unused trailing rows entirely.
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

Write-Host "Holdings comparison sheet refreshed."
